{"js": "// Remove the \"\u0399.\u039a.\u03a5.\" bullet item from the \u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397 (notification) list.\n// The whole paragraph (including its paragraph mark) is deleted so the\n// following list item (\"${local_directorate}\") simply moves up to take its\n// place, keeping its own (identical) list/paragraph formatting.\n\nconst body = context.document.body;\n\n// Locate the paragraph that contains the \"\u0399.\u039a.\u03a5.\" text.\nconst results = body.search(\"\u0399.\u039a.\u03a5.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const targetRange = results.items[0];\n  const paragraphs = targetRange.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  // Delete the whole paragraph (text + paragraph mark).\n  paragraphs.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"\u0399.\u039a.\u03a5.\" bullet item from the \u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397 (notification) list.\n# The whole paragraph (including its paragraph mark) is deleted so the\n# following list item (\"${local_directorate}\") simply moves up to take its\n# place, keeping its own (identical) list/paragraph formatting.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\u0399.\u039a.\u03a5.\", $true)\n\nif ($found) {\n    # Expand the found (collapsed) range to the whole paragraph, which in\n    # the Word object model includes the trailing paragraph mark, so the\n    # paragraph is fully removed (it does not just become empty).\n    $rng.Expand(4)  # wdParagraph\n    $rng.Delete()\n}\n"}
